$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 145, shifting existing rows 145:197 down to 146:198.
# This mirrors the other rows in the table (preserves the date number-format
# style that lives on column D).
$ws.Rows(145).Insert()

# Populate the newly inserted row 145 with the new observation
# (same shape as the surrounding "Primera" / China / $/caja 10 kilos rows).
$ws.Range("A145").Value = 11
$ws.Range("B145").Value = "Vega Monumental Concepción"
$ws.Range("C145").Value = "Bíobío"
$ws.Range("D145").Value = 44784
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = 100112003
$ws.Range("G145").Value = "Ajo"
$ws.Range("H145").Value = "Chino"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 400
$ws.Range("K145").Value = 23000
$ws.Range("L145").Value = 24000
$ws.Range("M145").Value = 23500
$ws.Range("N145").Value = "$/caja 10 kilos"
$ws.Range("O145").Value = "China"
$ws.Range("P145").Value = 2350
$ws.Range("Q145").Value = 10
$ws.Range("R145").Value = "Hortaliza"
